# HASSUYP-351: Korjaa migraatio
# Update the migration OID and status values in the template's second row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "1.2.246.578.5.1.2978288874.2711575506"
$ws.Range("B2").Value = "NAHTAVILLAOLO"
